$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efnb2"
$ws.Range("C2").Value = "Ephb6"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 20.94432133333333
$ws.Range("H2").Value = 62.832964
$ws.Range("I2").Value = 0.7396577289668299
$ws.Range("J2").Value = 0.7396577289668298
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.095455
$ws.Range("N2").Value = 0.286365
$ws.Range("O2").Value = 0.03201851307709132
$ws.Range("P2").Value = 0.03201851307709131
$ws.Range("Q2").Value = 1.999240192873333
$ws.Range("R2").Value = 17.99316173586
$ws.Range("S2").Value = 0.02368274066749611
$ws.Range("T2").Value = 0.0236827406674961

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efnb2"
$ws.Range("C3").Value = "Ephb6"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 20.94432133333333
$ws.Range("H3").Value = 62.832964
$ws.Range("I3").Value = 0.7396577289668299
$ws.Range("J3").Value = 0.7396577289668298
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.220310333333333
$ws.Range("N3").Value = 3.660931
$ws.Range("O3").Value = 0.4093292375039861
$ws.Range("P3").Value = 0.409329237503986
$ws.Range("Q3").Value = 25.55857174772044
$ws.Range("R3").Value = 230.027145729484
$ws.Range("S3").Value = 0.3027635342119225
$ws.Range("T3").Value = 0.3027635342119224

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efnb2"
$ws.Range("C4").Value = "Ephb6"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 20.94432133333333
$ws.Range("H4").Value = 62.832964
$ws.Range("I4").Value = 0.7396577289668299
$ws.Range("J4").Value = 0.7396577289668298
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.665478666666666
$ws.Range("N4").Value = 4.996435999999999
$ws.Range("O4").Value = 0.5586522494189227
$ws.Range("P4").Value = 0.5586522494189227
$ws.Range("Q4").Value = 34.88232036847822
$ws.Range("R4").Value = 313.940883316304
$ws.Range("S4").Value = 0.4132114540874113
$ws.Range("T4").Value = 0.4132114540874113

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efnb2"
$ws.Range("C5").Value = "Ephb6"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.327094666666667
$ws.Range("H5").Value = 6.981284
$ws.Range("I5").Value = 0.08218235047311259
$ws.Range("J5").Value = 0.08218235047311258
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.095455
$ws.Range("N5").Value = 0.286365
$ws.Range("O5").Value = 0.03201851307709132
$ws.Range("P5").Value = 0.03201851307709131
$ws.Range("Q5").Value = 0.2221328214066667
$ws.Range("R5").Value = 1.99919539266
$ws.Range("S5").Value = 0.002631356663329457
$ws.Range("T5").Value = 0.002631356663329457

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efnb2"
$ws.Range("C6").Value = "Ephb6"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.327094666666667
$ws.Range("H6").Value = 6.981284
$ws.Range("I6").Value = 0.08218235047311259
$ws.Range("J6").Value = 0.08218235047311258
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.220310333333333
$ws.Range("N6").Value = 3.660931
$ws.Range("O6").Value = 0.4093292375039861
$ws.Range("P6").Value = 0.409329237503986
$ws.Range("Q6").Value = 2.839777668378222
$ws.Range("R6").Value = 25.557999015404
$ws.Range("S6").Value = 0.03363963885544453
$ws.Range("T6").Value = 0.03363963885544452

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efnb2"
$ws.Range("C7").Value = "Ephb6"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.327094666666667
$ws.Range("H7").Value = 6.981284
$ws.Range("I7").Value = 0.08218235047311259
$ws.Range("J7").Value = 0.08218235047311258
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.665478666666666
$ws.Range("N7").Value = 4.996435999999999
$ws.Range("O7").Value = 0.5586522494189227
$ws.Range("P7").Value = 0.5586522494189227
$ws.Range("Q7").Value = 3.875726522647111
$ws.Range("R7").Value = 34.881538703824
$ws.Range("S7").Value = 0.04591135495433861
$ws.Range("T7").Value = 0.0459113549543386

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Efnb2"
$ws.Range("C8").Value = "Ephb6"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.044818
$ws.Range("H8").Value = 15.134454
$ws.Range("I8").Value = 0.1781599205600575
$ws.Range("J8").Value = 0.1781599205600575
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.095455
$ws.Range("N8").Value = 0.286365
$ws.Range("O8").Value = 0.03201851307709132
$ws.Range("P8").Value = 0.03201851307709131
$ws.Range("Q8").Value = 0.48155310219
$ws.Range("R8").Value = 4.333977919710001
$ws.Range("S8").Value = 0.005704415746265752
$ws.Range("T8").Value = 0.005704415746265751

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Efnb2"
$ws.Range("C9").Value = "Ephb6"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.044818
$ws.Range("H9").Value = 15.134454
$ws.Range("I9").Value = 0.1781599205600575
$ws.Range("J9").Value = 0.1781599205600575
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.220310333333333
$ws.Range("N9").Value = 3.660931
$ws.Range("O9").Value = 0.4093292375039861
$ws.Range("P9").Value = 0.409329237503986
$ws.Range("Q9").Value = 6.156243535185999
$ws.Range("R9").Value = 55.406191816674
$ws.Range("S9").Value = 0.07292606443661909
$ws.Range("T9").Value = 0.07292606443661907

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Efnb2"
$ws.Range("C10").Value = "Ephb6"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 5.044818
$ws.Range("H10").Value = 15.134454
$ws.Range("I10").Value = 0.1781599205600575
$ws.Range("J10").Value = 0.1781599205600575
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.665478666666666
$ws.Range("N10").Value = 4.996435999999999
$ws.Range("O10").Value = 0.5586522494189227
$ws.Range("P10").Value = 0.5586522494189227
$ws.Range("Q10").Value = 8.402036756215999
$ws.Range("R10").Value = 75.618330805944
$ws.Range("S10").Value = 0.0995294403771727
$ws.Range("T10").Value = 0.0995294403771727
